$wb = $excel.ActiveWorkbook

# Offense sheet ("OFF") - Row 3 corresponds to the "R" (road/running) totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 173
$wsOff.Range("C3").Value = 129
$wsOff.Range("D3").Value = 47
$wsOff.Range("E3").Value = 21
$wsOff.Range("G3").Value = 6

# Defense sheet ("DEF") - Row 3 corresponds to the "R" (road/running) totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 477
$wsDef.Range("C3").Value = 343
$wsDef.Range("D3").Value = 111
$wsDef.Range("E3").Value = 52
